$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "jairidssas@yahoo.com"
$ws.Range("D2").Select()
